$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A23").Value = 83
$ws.Range("D9").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("C23").Value = "יש שריפות ענק בצפון. בתחקיר ראשון עולה שזאת הצתה או מנגל"
$f = $ws.Range("C23").Font
$f.Color = 0
$f.Size = 15
$ws.Rows.Item(23).AutoFit()
$h = $ws.Rows.Item(23).RowHeight
Write-Host ("autofit height=" + $h)
